$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.266.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.877.40"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.09%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5049"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3936"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.88%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09564"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.25%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.144"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.02"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.496"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.98"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.99%  "
$ws.Range("B15").Value = "BinanceUSD"
$ws.Range("C15").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.005"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.439"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.76%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.40%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06592"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.197"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.326.71"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.46%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.299"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.571"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.095.48"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.29"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +5.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "158.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.53"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.44%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.070"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.648"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.620"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.555"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06728"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("E37").Value = "  +4.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2198"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.21%  "
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.50"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.18%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6363"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.007"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.185"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.000"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.48"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5993"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.658"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.270"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.006"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.196"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.03%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06870"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.72%  "
